$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-strings text order puts "ECs" before "FAPs" (ECs, FAPs, sCs, Sema3c, Plxnd1).
# Build the full target table: sending cluster in {ECs, FAPs, sCs} x Ligand=Sema3c x
# Receptor=Plxnd1 x target cluster in {ECs, FAPs, sCs}, 9 data rows total (rows 2-10).

$ligand = "Sema3c"
$receptor = "Plxnd1"
$clusters = @("ECs", "FAPs", "sCs")

# E .. T numeric values per (sending cluster, target cluster) combination.
$rowData = @{
  "ECs|ECs"   = @(1, 0.3333333333333333, 0.2364093333333333, 0.709228, 0.005805733041453686, 0.005805733041453687, 3, 1, 88.37814633333333, 265.134439, 0.7138016014383547, 0.7138016014383547, 20.89341865589911, 188.040767903092, 0.004144141542513211, 0.004144141542513212)
  "ECs|FAPs"  = @(1, 0.3333333333333333, 0.2364093333333333, 0.709228, 0.005805733041453686, 0.005805733041453687, 3, 1, 8.866675333333333, 26.600026, 0.07161325864989523, 0.07161325864989525, 2.096164804436444, 18.865483239928, 0.0004157674619498658, 0.0004157674619498659)
  "ECs|sCs"   = @(1, 0.3333333333333333, 0.2364093333333333, 0.709228, 0.005805733041453686, 0.005805733041453687, 3, 1, 26.56849866666667, 79.705496, 0.2145851399117501, 0.2145851399117501, 6.281041057454222, 56.52936951708799, 0.00124582403699061, 0.00124582403699061)
  "FAPs|ECs"  = @(3, 1, 39.09670133333334, 117.290104, 0.9601355730855794, 0.9601355730855795, 3, 1, 88.37814633333333, 265.134439, 0.7138016014383547, 0.7138016014383547, 3455.293991587962, 31097.64592429166, 0.6853463096664191, 0.6853463096664191)
  "FAPs|FAPs" = @(3, 1, 39.09670133333334, 117.290104, 0.9601355730855794, 0.9601355730855795, 3, 1, 8.866675333333333, 26.600026, 0.07161325864989523, 0.07161325864989525, 346.6577573269671, 3119.919815942704, 0.06875843713434299, 0.068758437134343)
  "FAPs|sCs"  = @(3, 1, 39.09670133333334, 117.290104, 0.9601355730855794, 0.9601355730855795, 3, 1, 26.56849866666667, 79.705496, 0.2145851399117501, 0.2145851399117501, 1038.740657245732, 9348.665915211584, 0.2060308262848174, 0.2060308262848174)
  "sCs|ECs"   = @(3, 1, 1.386869333333333, 4.160608, 0.03405869387296686, 0.03405869387296687, 3, 1, 88.37814633333333, 265.134439, 0.7138016014383547, 0.7138016014383547, 122.5689408865458, 1103.120467978912, 0.02431115022942242, 0.02431115022942243)
  "sCs|FAPs"  = @(3, 1, 1.386869333333333, 4.160608, 0.03405869387296686, 0.03405869387296687, 3, 1, 8.866675333333333, 26.600026, 0.07161325864989523, 0.07161325864989525, 12.29692010842311, 110.672280975808, 0.002439054053602378, 0.002439054053602379)
  "sCs|sCs"   = @(3, 1, 1.386869333333333, 4.160608, 0.03405869387296686, 0.03405869387296687, 3, 1, 26.56849866666667, 79.705496, 0.2145851399117501, 0.2145851399117501, 36.84703603350756, 331.6233243015679, 0.007308489589942059, 0.007308489589942061)
}

$r = 2
foreach ($sending in $clusters) {
  foreach ($target in $clusters) {
    $ws.Cells.Item($r, 1).Value = $sending
    $ws.Cells.Item($r, 2).Value = $ligand
    $ws.Cells.Item($r, 3).Value = $receptor
    $ws.Cells.Item($r, 4).Value = $target

    $vals = $rowData["$sending|$target"]
    for ($i = 0; $i -lt $vals.Length; $i++) {
      $ws.Cells.Item($r, 5 + $i).Value = $vals[$i]
    }

    $r++
  }
}
